$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The enemy-stat columns are being reordered: DetectionRange now comes
# before RotationSpeed, and a brand-new "AttackRotationSpeed" AI column is
# inserted right after RotationSpeed (so enemies can rotate toward the
# player while attacking instead of orbiting them).

# 1) Swap the RotationSpeed (E) and DetectionRange (F) columns - cut F and
#    insert it in front of E. This carries the header, the per-row values
#    and the column width along with the move.
$ws.Range("F1").EntireColumn.Cut() | Out-Null
$ws.Range("E1").EntireColumn.Insert() | Out-Null

# 2) Insert the new AttackRotationSpeed column after RotationSpeed (now
#    column F), pushing AttackRate..DeathAnimName one column to the right.
$ws.Range("G1").EntireColumn.Insert() | Out-Null

$ws.Range("G1").Value2 = "AttackRotationSpeed"
$ws.Range("G2").Value2 = 250
$ws.Range("G3").Value2 = 250
$ws.Range("G4").Value2 = 10

# Match the new column's width to its neighbours (~20.375 chars wide).
$ws.Range("G1:G4").EntireColumn.ColumnWidth = 19.7

# Reflect the editor's final selection.
$ws.Range("G5").Select()
